$wb = $excel.ActiveWorkbook

$wsMenu = $wb.Worksheets.Item("Menu Mock")
$wsProviders = $wb.Worksheets.Item("Providers")

# --- Menu Mock sheet updates (columns D and E) ---
$wsMenu.Range("D42").Value = "0.72s"
$wsMenu.Range("E42").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D43").Value = "3.5s"
$wsMenu.Range("E43").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D44").Value = "3.6x"
$wsMenu.Range("E44").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D50").Value = "1.125s"
$wsMenu.Range("E50").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D51").Value = "3.5s"
$wsMenu.Range("E51").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D52").Value = "3.6x"
$wsMenu.Range("E52").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D58").Value = "1.125s"
$wsMenu.Range("E58").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D59").Value = "3.5s"
$wsMenu.Range("E59").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D60").Value = "3.6x"
$wsMenu.Range("E60").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D66").Value = "1.4s"
$wsMenu.Range("E66").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D67").Value = "2.8s"
$wsMenu.Range("E67").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D68").Value = "3x"
$wsMenu.Range("E68").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D74").Value = "1.8s"
$wsMenu.Range("E74").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("E75").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D76").Value = "2x"
$wsMenu.Range("E76").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D81").Value = "2.1s"
$wsMenu.Range("E81").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D82").Value = "7.2s"
$wsMenu.Range("E82").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D83").Value = "2x"
$wsMenu.Range("E83").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
$wsMenu.Range("D88").Value = "0.90s"
$wsMenu.Range("E88").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsMenu.Range("D89").Value = "4.2s"
$wsMenu.Range("E89").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsMenu.Range("D90").Value = "4x"
$wsMenu.Range("E90").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"

# --- Providers sheet updates (column B) ---
$wsProviders.Range("B10").Value = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 19.6s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"
$wsProviders.Range("B11").Value = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.8s | 3.0s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"
$wsProviders.Range("B12").Value = "1.6x | 1.8x | 2x | 2.4x | 2.7x | 3x | 3.2x | 3.6x | 4x | 4.5x | 5x | 6x | 7.5x | 8x | 9x | 10x | 12x | 12.5x | 15x | 16x | 20x | 25x"
